# Generate Report for Archive
#
# The localization-status report is regenerated: the handoff record that
# was still showing "Ready for handoff" has moved on to "In Translation".
# That status string is shared by the Overview rollup sheet (columns
# zh-cn/de-de) as well as each per-language detail sheet's "Status"
# column, so update it everywhere it appears. The shorter replacement
# text lets Excel shrink those "Status" columns, so resize them to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The "Status" columns narrow to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
